$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.915919409345647
$ws.Range("D2").Value = 9.169527626918649
$ws.Range("E2").Value = 13.68066547526212
$ws.Range("F2").Value = 33.0841991395868
$ws.Range("G2").Value = 3.657129963734116
$ws.Range("J2").Value = 9.941241086588521
$ws.Range("K2").Value = 17.73700119242937
$ws.Range("O2").Value = 24.95574771454122
$ws.Range("C3").Value = 4.750707899966622
$ws.Range("D3").Value = 9.092012421453571
$ws.Range("E3").Value = 13.61042338607188
$ws.Range("F3").Value = 33.21547950151452
$ws.Range("G3").Value = 3.659708854907774
$ws.Range("J3").Value = 9.951379488853087
$ws.Range("K3").Value = 17.08162149285088
$ws.Range("O3").Value = 25.10830927973124
$ws.Range("C4").Value = 4.647826968361152
$ws.Range("D4").Value = 9.045122793239628
$ws.Range("E4").Value = 13.56977425430711
$ws.Range("F4").Value = 33.30735863069133
$ws.Range("G4").Value = 3.661373982513858
$ws.Range("J4").Value = 9.959316335422155
$ws.Range("K4").Value = 16.66660771826491
$ws.Range("O4").Value = 25.20997124965395
$ws.Range("C5").Value = 4.605615247610249
$ws.Range("D5").Value = 9.026207434453118
$ws.Range("E5").Value = 13.55384539326437
$ws.Range("F5").Value = 33.34761792905387
$ws.Range("G5").Value = 3.662073143693813
$ws.Range("J5").Value = 9.962980967987116
$ws.Range("K5").Value = 16.49454223247431
$ws.Range("O5").Value = 25.25339847164798
$ws.Range("C6").Value = 4.598591005012123
$ws.Range("D6").Value = 9.023078622180794
$ws.Range("E6").Value = 13.55123917612368
$ws.Range("F6").Value = 33.35447264213537
$ws.Range("G6").Value = 3.662190485601744
$ws.Range("J6").Value = 9.963615460262547
$ws.Range("K6").Value = 16.46580017710194
$ws.Range("O6").Value = 25.26072999165884
$ws.Range("C7").Value = 4.647258744366158
$ws.Range("D7").Value = 9.044866894563441
$ws.Range("E7").Value = 13.56955684124521
$ws.Range("F7").Value = 33.30789019514409
$ws.Range("G7").Value = 3.661383328112728
$ws.Range("J7").Value = 9.959364015936201
$ws.Range("K7").Value = 16.66429878596729
$ws.Range("O7").Value = 25.21054884266686
$ws.Range("C8").Value = 4.859296850483752
$ws.Range("D8").Value = 9.142663307536976
$ws.Range("E8").Value = 13.65593903580909
$ws.Range("F8").Value = 33.12711512788927
$ws.Range("G8").Value = 3.658002255119725
$ws.Range("J8").Value = 9.944381456768554
$ws.Range("K8").Value = 17.51377135775846
$ws.Range("O8").Value = 25.00668710350229
$ws.Range("C9").Value = 5.260659255440519
$ws.Range("D9").Value = 9.339352729403515
$ws.Range("E9").Value = 13.84442591098164
$ws.Range("F9").Value = 32.8627774812246
$ws.Range("G9").Value = 3.652016876501484
$ws.Range("J9").Value = 9.928586846664809
$ws.Range("K9").Value = 19.07081393553876
$ws.Range("O9").Value = 24.67073398291707
$ws.Range("C10").Value = 5.543117735967822
$ws.Range("D10").Value = 9.485926335229905
$ws.Range("E10").Value = 13.99371153109783
$ws.Range("F10").Value = 32.72447292377322
$ws.Range("G10").Value = 3.64800809603376
$ws.Range("J10").Value = 9.925265249519514
$ws.Range("K10").Value = 20.13835958286496
$ws.Range("O10").Value = 24.46341203797078
$ws.Range("C11").Value = 5.668271961069204
$ws.Range("D11").Value = 9.552855370451514
$ws.Range("E11").Value = 14.06377916599529
$ws.Range("F11").Value = 32.67387858525831
$ws.Range("G11").Value = 3.646267842190086
$ws.Range("J11").Value = 9.925550467838756
$ws.Range("K11").Value = 20.6057005986726
$ws.Range("O11").Value = 24.3778031677684
$ws.Range("C12").Value = 5.715138505885783
$ws.Range("D12").Value = 9.57821905980531
$ws.Range("E12").Value = 14.09060496738707
$ws.Range("F12").Value = 32.65650539422367
$ws.Range("G12").Value = 3.645620767995213
$ws.Range("J12").Value = 9.925916343610194
$ws.Range("K12").Value = 20.77992125378159
$ws.Range("O12").Value = 24.34664767425145
$ws.Range("C13").Value = 5.705069066962375
$ws.Range("D13").Value = 9.572755956686605
$ws.Range("E13").Value = 14.08481481575103
$ws.Range("F13").Value = 32.66016738542265
$ws.Range("G13").Value = 3.645759597789268
$ws.Range("J13").Value = 9.925826085008902
$ws.Range("K13").Value = 20.74252383522916
$ws.Range("O13").Value = 24.35330122731838
$ws.Range("C14").Value = 5.6721385033209
$ws.Range("D14").Value = 9.554941770987764
$ws.Range("E14").Value = 14.06598037616885
$ws.Range("F14").Value = 32.6724134156669
$ws.Range("G14").Value = 3.646214368439254
$ws.Range("J14").Value = 9.92557540345036
$ws.Range("K14").Value = 20.6200896111809
$ws.Range("O14").Value = 24.37521461498207
$ws.Range("C15").Value = 5.651897702231847
$ws.Range("D15").Value = 9.544032051126534
$ws.Range("E15").Value = 14.05448133122287
$ws.Range("F15").Value = 32.6801474246742
$ws.Range("G15").Value = 3.646494479290009
$ws.Range("J15").Value = 9.925455421229916
$ws.Range("K15").Value = 20.54473350769399
$ws.Range("O15").Value = 24.38880197696864
$ws.Range("C16").Value = 5.534867337760269
$ws.Range("D16").Value = 9.481555970360837
$ws.Range("E16").Value = 13.9891743682202
$ws.Range("F16").Value = 32.72802853595638
$ws.Range("G16").Value = 3.648123497577645
$ws.Range("J16").Value = 9.925282721030358
$ws.Range("K16").Value = 20.10743896389541
$ws.Range("O16").Value = 24.46918286529985
$ws.Range("C17").Value = 5.462183225966606
$ws.Range("D17").Value = 9.443281248577463
$ws.Range("E17").Value = 13.94965083689505
$ws.Range("F17").Value = 32.76056806097219
$ws.Range("G17").Value = 3.649144151807302
$ws.Range("J17").Value = 9.925636569398202
$ws.Range("K17").Value = 19.83439675584778
$ws.Range("O17").Value = 24.52073095871918
$ws.Range("C18").Value = 5.420065540779516
$ws.Range("D18").Value = 9.421291377964867
$ws.Range("E18").Value = 13.92712203961552
$ws.Range("F18").Value = 32.78044287515694
$ws.Range("G18").Value = 3.649739055267031
$ws.Range("J18").Value = 9.926009196649222
$ws.Range("K18").Value = 19.67563459926921
$ws.Range("O18").Value = 24.55119854197805
$ws.Range("C19").Value = 5.405753191677012
$ws.Range("D19").Value = 9.413850757159087
$ws.Range("E19").Value = 13.91952975217727
$ws.Range("F19").Value = 32.78737075174886
$ws.Range("G19").Value = 3.649941829613085
$ws.Range("J19").Value = 9.926164416364362
$ws.Range("K19").Value = 19.62158978348779
$ws.Range("O19").Value = 24.56165457531961
$ws.Range("C20").Value = 5.469953189532563
$ws.Range("D20").Value = 9.447353222158382
$ws.Range("E20").Value = 13.95383718623765
$ws.Range("F20").Value = 32.75698412227514
$ws.Range("G20").Value = 3.649034689433297
$ws.Range("J20").Value = 9.925581402790048
$ws.Range("K20").Value = 19.86364104197167
$ws.Range("O20").Value = 24.5151587847488
$ws.Range("C21").Value = 5.681825647008546
$ws.Range("D21").Value = 9.560173846229015
$ws.Range("E21").Value = 14.07150470506023
$ws.Range("F21").Value = 32.66876788702755
$ws.Range("G21").Value = 3.646080468227802
$ws.Range("J21").Value = 9.925642040206659
$ws.Range("K21").Value = 20.65612706236585
$ws.Range("O21").Value = 24.36874376182269
$ws.Range("C22").Value = 5.817208867736066
$ws.Range("D22").Value = 9.634012219045671
$ws.Range("E22").Value = 14.15010494002865
$ws.Range("F22").Value = 32.6215282176645
$ws.Range("G22").Value = 3.64421917552344
$ws.Range("J22").Value = 9.927184446087553
$ws.Range("K22").Value = 21.1579811639151
$ws.Range("O22").Value = 24.28041847591739
$ws.Range("C23").Value = 5.745248893110666
$ws.Range("D23").Value = 9.594599437195328
$ws.Range("E23").Value = 14.1080050609977
$ws.Range("F23").Value = 32.6457835892232
$ws.Range("G23").Value = 3.645206248094084
$ws.Range("J23").Value = 9.926223904209127
$ws.Range("K23").Value = 20.89163947123648
$ws.Range("O23").Value = 24.32688167874613
$ws.Range("C24").Value = 5.466441415258419
$ws.Range("D24").Value = 9.445512235596629
$ws.Range("E24").Value = 13.95194393271576
$ws.Range("F24").Value = 32.75860078551246
$ws.Range("G24").Value = 3.649084152086479
$ws.Range("J24").Value = 9.925605816552006
$ws.Range("K24").Value = 19.85042525559416
$ws.Range("O24").Value = 24.5176753740212
$ws.Range("C25").Value = 5.154023060627293
$ws.Range("D25").Value = 9.285712484793837
$ws.Range("E25").Value = 13.79147453579053
$ws.Range("F25").Value = 32.92453538994147
$ws.Range("G25").Value = 3.653567502841268
$ws.Range("J25").Value = 9.931404732119876
$ws.Range("K25").Value = 18.66236352704335
$ws.Range("O25").Value = 24.754724580055
